$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data range (before the existing
# row 47), pushing the existing rows 47-143 down to 49-145.
$ws.Rows("47:48").Insert()

# --- New row 47 ---
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value = "Arica y Parinacota"
$ws.Range("D47").Value = 45114
$ws.Range("E47").Value = 15
$ws.Range("F47").Value = 100112021
$ws.Range("G47").Value = "Ají"
$ws.Range("H47").Value = "Inferno"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 140
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 11000
$ws.Range("M47").Value = 10500
$ws.Range("N47").Value = "$/caja 15 kilos"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value = 700
$ws.Range("Q47").Value = 15
$ws.Range("R47").Value = "Hortaliza"

# --- New row 48 ---
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 45114
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100112021
$ws.Range("G48").Value = "Ají"
$ws.Range("H48").Value = "Inferno"
$ws.Range("I48").Value = "Segunda"
$ws.Range("J48").Value = 170
$ws.Range("K48").Value = 8000
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = 8500
$ws.Range("N48").Value = "$/caja 15 kilos"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 567
$ws.Range("Q48").Value = 15
$ws.Range("R48").Value = "Hortaliza"
